# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows for "Piña" (Vega Modelo de Temuco)
# right before the existing row 346, shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 346 (pushes old rows 346.. down to 348..)
$ws.Rows.Item(346).Insert()
$ws.Rows.Item(346).Insert()

# --- New row 346: Primera quality ---
$ws.Cells.Item(346, 1).Value = 10
$ws.Cells.Item(346, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(346, 3).Value = "La Araucanía"
$ws.Cells.Item(346, 4).Value = 44642
$ws.Cells.Item(346, 5).Value = 9
$ws.Cells.Item(346, 6).Value = "Fruta"
$ws.Cells.Item(346, 7).Value = 100108
$ws.Cells.Item(346, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(346, 9).Value = 100108005
$ws.Cells.Item(346, 10).Value = "Piña"
$ws.Cells.Item(346, 11).Value = "Caramelo"
$ws.Cells.Item(346, 12).Value = "Primera"
$ws.Cells.Item(346, 13).Value = 100
$ws.Cells.Item(346, 14).Value = 18000
$ws.Cells.Item(346, 15).Value = 18000
$ws.Cells.Item(346, 16).Value = 18000
$ws.Cells.Item(346, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(346, 18).Value = "Ecuador"
$ws.Cells.Item(346, 19).Value = 1500
$ws.Cells.Item(346, 20).Value = 12

# --- New row 347: Segunda quality ---
$ws.Cells.Item(347, 1).Value = 10
$ws.Cells.Item(347, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(347, 3).Value = "La Araucanía"
$ws.Cells.Item(347, 4).Value = 44642
$ws.Cells.Item(347, 5).Value = 9
$ws.Cells.Item(347, 6).Value = "Fruta"
$ws.Cells.Item(347, 7).Value = 100108
$ws.Cells.Item(347, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(347, 9).Value = 100108005
$ws.Cells.Item(347, 10).Value = "Piña"
$ws.Cells.Item(347, 11).Value = "Caramelo"
$ws.Cells.Item(347, 12).Value = "Segunda"
$ws.Cells.Item(347, 13).Value = 50
$ws.Cells.Item(347, 14).Value = 18000
$ws.Cells.Item(347, 15).Value = 18000
$ws.Cells.Item(347, 16).Value = 18000
$ws.Cells.Item(347, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(347, 18).Value = "Ecuador"
$ws.Cells.Item(347, 19).Value = 1286
$ws.Cells.Item(347, 20).Value = 14
